# Auto-generated edit script applying the Famfrit_Profits data refresh diff.
# Updates cached market-board values (H..N) across all 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 4735.174
$ws.Cells.Item(9, 9).Value = 5445.2104
$ws.Cells.Item(9, 10).Value = 1362.5
$ws.Cells.Item(9, 11).Value = 5445.2104
$ws.Cells.Item(9, 12).Value = 1362.5
$ws.Cells.Item(9, 13).Value = -5276.2104
$ws.Cells.Item(9, 14).Value = -1700.5
$ws.Cells.Item(17, 8).Value = 1731393.1
$ws.Cells.Item(17, 10).Value = 1731393.1
$ws.Cells.Item(17, 12).Value = 5194179.300000001
$ws.Cells.Item(17, 14).Value = -5194515.300000001
$ws.Cells.Item(47, 8).Value = 17667
$ws.Cells.Item(47, 9).Value = 17667
$ws.Cells.Item(47, 11).Value = 17667
$ws.Cells.Item(47, 13).Value = -16695
$ws.Cells.Item(51, 8).Value = 3961.0356
$ws.Cells.Item(51, 9).Value = 2401.6365
$ws.Cells.Item(51, 10).Value = 4970.0586
$ws.Cells.Item(51, 11).Value = 2401.6365
$ws.Cells.Item(51, 12).Value = 4970.0586
$ws.Cells.Item(51, 13).Value = -1917.6365
$ws.Cells.Item(51, 14).Value = -5938.0586
$ws.Cells.Item(116, 8).Value = 5899.8335
$ws.Cells.Item(116, 9).Value = 6075
$ws.Cells.Item(116, 10).Value = 5549.5
$ws.Cells.Item(116, 11).Value = 6075
$ws.Cells.Item(116, 12).Value = 5549.5
$ws.Cells.Item(116, 13).Value = -2633
$ws.Cells.Item(116, 14).Value = -12433.5
$ws.Cells.Item(132, 8).Value = 3724.923
$ws.Cells.Item(132, 9).Value = 3627
$ws.Cells.Item(132, 11).Value = 10881
$ws.Cells.Item(132, 13).Value = -8351
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(133, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 130000
$ws.Cells.Item(134, 10).Value = 130000
$ws.Cells.Item(134, 12).Value = 130000
$ws.Cells.Item(134, 14).Value = -140140
$ws.Cells.Item(135, 8).Value = 2845.6
$ws.Cells.Item(135, 9).Value = 1173.8334
$ws.Cells.Item(135, 11).Value = 10564.5006
$ws.Cells.Item(135, 13).Value = -8029.500599999999
$ws.Cells.Item(137, 8).Value = 3639.88
$ws.Cells.Item(137, 9).Value = 2004.7273
$ws.Cells.Item(137, 10).Value = 4924.643
$ws.Cells.Item(137, 11).Value = 6014.1819
$ws.Cells.Item(137, 12).Value = 14773.929
$ws.Cells.Item(137, 13).Value = -3464.1819
$ws.Cells.Item(137, 14).Value = -19873.929
$ws.Cells.Item(138, 8).Value = 6768.64
$ws.Cells.Item(138, 9).Value = 2627.3
$ws.Cells.Item(138, 10).Value = 7803.975
$ws.Cells.Item(138, 11).Value = 7881.900000000001
$ws.Cells.Item(138, 12).Value = 23411.925
$ws.Cells.Item(138, 13).Value = -2741.900000000001
$ws.Cells.Item(138, 14).Value = -33691.925
$ws.Cells.Item(140, 8).Value = 116347.5
$ws.Cells.Item(140, 9).Value = 65000
$ws.Cells.Item(140, 10).Value = 133463.33
$ws.Cells.Item(140, 11).Value = 65000
$ws.Cells.Item(140, 12).Value = 133463.33
$ws.Cells.Item(140, 13).Value = -59820
$ws.Cells.Item(140, 14).Value = -143823.33
$ws.Cells.Item(141, 8).Value = 5089.9165
$ws.Cells.Item(141, 9).Value = 6242.269
$ws.Cells.Item(141, 10).Value = 2093.8
$ws.Cells.Item(141, 11).Value = 18726.807
$ws.Cells.Item(141, 12).Value = 6281.400000000001
$ws.Cells.Item(141, 13).Value = -13546.807
$ws.Cells.Item(141, 14).Value = -16641.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2438.1
$ws.Cells.Item(2, 9).Value = 1727.3077
$ws.Cells.Item(2, 10).Value = 3758.1428
$ws.Cells.Item(2, 11).Value = 1727.3077
$ws.Cells.Item(2, 12).Value = 3758.1428
$ws.Cells.Item(2, 13).Value = -1614.3077
$ws.Cells.Item(2, 14).Value = -3984.1428
$ws.Cells.Item(32, 8).Value = 2267.3538
$ws.Cells.Item(32, 9).Value = 1368.9178
$ws.Cells.Item(32, 10).Value = 9554.666999999999
$ws.Cells.Item(32, 11).Value = 1368.9178
$ws.Cells.Item(32, 12).Value = 9554.666999999999
$ws.Cells.Item(32, 13).Value = -1081.9178
$ws.Cells.Item(32, 14).Value = -10128.667
$ws.Cells.Item(37, 8).Value = 33268
$ws.Cells.Item(37, 9).Value = 30599.2
$ws.Cells.Item(37, 10).Value = 39940
$ws.Cells.Item(37, 11).Value = 30599.2
$ws.Cells.Item(37, 12).Value = 39940
$ws.Cells.Item(37, 13).Value = -30326.2
$ws.Cells.Item(37, 14).Value = -40486
$ws.Cells.Item(45, 8).Value = 3890
$ws.Cells.Item(45, 10).Value = 1794.6666
$ws.Cells.Item(45, 12).Value = 1794.6666
$ws.Cells.Item(45, 14).Value = -2548.6666
$ws.Cells.Item(61, 8).Value = 17858810
$ws.Cells.Item(61, 9).Value = 20834846
$ws.Cells.Item(61, 11).Value = 20834846
$ws.Cells.Item(61, 13).Value = -20834634
$ws.Cells.Item(74, 8).Value = 21301732
$ws.Cells.Item(74, 9).Value = 24418656
$ws.Cells.Item(74, 10).Value = 2740.1667
$ws.Cells.Item(74, 11).Value = 24418656
$ws.Cells.Item(74, 12).Value = 2740.1667
$ws.Cells.Item(74, 13).Value = -24417782
$ws.Cells.Item(74, 14).Value = -4488.1667
$ws.Cells.Item(77, 8).Value = 21301732
$ws.Cells.Item(77, 9).Value = 24418656
$ws.Cells.Item(77, 10).Value = 2740.1667
$ws.Cells.Item(77, 11).Value = 122093280
$ws.Cells.Item(77, 12).Value = 13700.8335
$ws.Cells.Item(77, 13).Value = -122088912
$ws.Cells.Item(77, 14).Value = -22436.8335
$ws.Cells.Item(102, 8).Value = 102060.55
$ws.Cells.Item(102, 9).Value = 106905.9
$ws.Cells.Item(102, 11).Value = 106905.9
$ws.Cells.Item(102, 13).Value = -105283.9
$ws.Cells.Item(116, 8).Value = 2438.1
$ws.Cells.Item(116, 9).Value = 1727.3077
$ws.Cells.Item(116, 10).Value = 3758.1428
$ws.Cells.Item(116, 11).Value = 1727.3077
$ws.Cells.Item(116, 12).Value = 3758.1428
$ws.Cells.Item(116, 13).Value = 566.6922999999999
$ws.Cells.Item(116, 14).Value = -8346.1428
$ws.Cells.Item(134, 8).Value = 100000
$ws.Cells.Item(134, 10).Value = 100000
$ws.Cells.Item(134, 12).Value = 100000
$ws.Cells.Item(134, 14).Value = -110140
$ws.Cells.Item(135, 8).Value = 71834.17999999999
$ws.Cells.Item(135, 10).Value = 71834.17999999999
$ws.Cells.Item(135, 12).Value = 71834.17999999999
$ws.Cells.Item(135, 14).Value = -81974.17999999999
$ws.Cells.Item(136, 8).Value = 17858810
$ws.Cells.Item(136, 9).Value = 20834846
$ws.Cells.Item(136, 11).Value = 62504538
$ws.Cells.Item(136, 13).Value = -62501988
$ws.Cells.Item(137, 8).Value = 100000
$ws.Cells.Item(137, 9).Value = 100000
$ws.Cells.Item(137, 11).Value = 100000
$ws.Cells.Item(137, 13).Value = -94900
$ws.Cells.Item(140, 8).Value = 100000
$ws.Cells.Item(140, 10).Value = 100000
$ws.Cells.Item(140, 12).Value = 100000
$ws.Cells.Item(140, 14).Value = -110360
$ws.Cells.Item(141, 8).Value = 100000
$ws.Cells.Item(141, 10).Value = 100000
$ws.Cells.Item(141, 12).Value = 100000
$ws.Cells.Item(141, 14).Value = -110360

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2438.1
$ws.Cells.Item(3, 9).Value = 1727.3077
$ws.Cells.Item(3, 10).Value = 3758.1428
$ws.Cells.Item(3, 11).Value = 1727.3077
$ws.Cells.Item(3, 12).Value = 3758.1428
$ws.Cells.Item(3, 13).Value = -1613.3077
$ws.Cells.Item(3, 14).Value = -3986.1428
$ws.Cells.Item(20, 8).Value = 1928.174
$ws.Cells.Item(20, 10).Value = 3009.4443
$ws.Cells.Item(20, 12).Value = 3009.4443
$ws.Cells.Item(20, 14).Value = -3503.4443
$ws.Cells.Item(86, 8).Value = 11055.333
$ws.Cells.Item(86, 9).Value = 15106.375
$ws.Cells.Item(86, 10).Value = 2953.25
$ws.Cells.Item(86, 11).Value = 15106.375
$ws.Cells.Item(86, 12).Value = 2953.25
$ws.Cells.Item(86, 13).Value = -13983.375
$ws.Cells.Item(86, 14).Value = -5199.25
$ws.Cells.Item(89, 8).Value = 11055.333
$ws.Cells.Item(89, 9).Value = 15106.375
$ws.Cells.Item(89, 10).Value = 2953.25
$ws.Cells.Item(89, 11).Value = 75531.875
$ws.Cells.Item(89, 12).Value = 14766.25
$ws.Cells.Item(89, 13).Value = -69915.875
$ws.Cells.Item(89, 14).Value = -25998.25
$ws.Cells.Item(105, 8).Value = 5032.6
$ws.Cells.Item(105, 9).Value = 6474.773
$ws.Cells.Item(105, 11).Value = 6474.773
$ws.Cells.Item(105, 13).Value = -4727.773
$ws.Cells.Item(132, 8).Value = 128264.25
$ws.Cells.Item(132, 10).Value = 128264.25
$ws.Cells.Item(132, 12).Value = 128264.25
$ws.Cells.Item(132, 14).Value = -138384.25
$ws.Cells.Item(133, 8).Value = 100585
$ws.Cells.Item(133, 10).Value = 100780
$ws.Cells.Item(133, 12).Value = 100780
$ws.Cells.Item(133, 14).Value = -110900
$ws.Cells.Item(134, 8).Value = 1886.6857
$ws.Cells.Item(134, 9).Value = 1654.5
$ws.Cells.Item(134, 10).Value = 3279.8
$ws.Cells.Item(134, 11).Value = 4963.5
$ws.Cells.Item(134, 12).Value = 9839.400000000001
$ws.Cells.Item(134, 13).Value = -2428.5
$ws.Cells.Item(134, 14).Value = -14909.4
$ws.Cells.Item(138, 8).Value = 196666
$ws.Cells.Item(138, 10).Value = 196666
$ws.Cells.Item(138, 12).Value = 196666
$ws.Cells.Item(138, 14).Value = -206946
$ws.Cells.Item(141, 8).Value = 119975
$ws.Cells.Item(141, 10).Value = 119975
$ws.Cells.Item(141, 12).Value = 119975
$ws.Cells.Item(141, 14).Value = -130335

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 9).Value = 188.6
$ws.Cells.Item(7, 10).Value = 416.42856
$ws.Cells.Item(7, 11).Value = 188.6
$ws.Cells.Item(7, 12).Value = 416.42856
$ws.Cells.Item(7, 13).Value = -75.59999999999999
$ws.Cells.Item(7, 14).Value = -642.4285600000001
$ws.Cells.Item(16, 8).Value = 2460.889
$ws.Cells.Item(16, 9).Value = 2457.4
$ws.Cells.Item(16, 10).Value = 2465.25
$ws.Cells.Item(16, 11).Value = 2457.4
$ws.Cells.Item(16, 12).Value = 2465.25
$ws.Cells.Item(16, 13).Value = -2170.4
$ws.Cells.Item(16, 14).Value = -3039.25
$ws.Cells.Item(62, 8).Value = 3000
$ws.Cells.Item(62, 9).Value = 3000
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 3000
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -2376
$ws.Cells.Item(62, 14).ClearContents()
$ws.Cells.Item(65, 8).Value = 3000
$ws.Cells.Item(65, 9).Value = 3000
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 15000
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -11880
$ws.Cells.Item(65, 14).ClearContents()
$ws.Cells.Item(76, 8).Value = 3040
$ws.Cells.Item(76, 9).Value = 3040
$ws.Cells.Item(76, 11).Value = 3040
$ws.Cells.Item(76, 13).Value = -2725
$ws.Cells.Item(79, 8).Value = 3040
$ws.Cells.Item(79, 9).Value = 3040
$ws.Cells.Item(79, 11).Value = 3040
$ws.Cells.Item(79, 13).Value = -1948
$ws.Cells.Item(113, 8).Value = 2460.889
$ws.Cells.Item(113, 9).Value = 2457.4
$ws.Cells.Item(113, 10).Value = 2465.25
$ws.Cells.Item(113, 11).Value = 2457.4
$ws.Cells.Item(113, 12).Value = 2465.25
$ws.Cells.Item(113, 13).Value = -287.4000000000001
$ws.Cells.Item(113, 14).Value = -6805.25
$ws.Cells.Item(132, 8).Value = 56249.297
$ws.Cells.Item(132, 9).Value = 68774.3
$ws.Cells.Item(132, 10).Value = 2570.7144
$ws.Cells.Item(132, 11).Value = 206322.9
$ws.Cells.Item(132, 12).Value = 7712.1432
$ws.Cells.Item(132, 13).Value = -203792.9
$ws.Cells.Item(132, 14).Value = -12772.1432
$ws.Cells.Item(133, 8).Value = 60162.5
$ws.Cells.Item(133, 9).Value = 35000
$ws.Cells.Item(133, 10).Value = 68550
$ws.Cells.Item(133, 11).Value = 35000
$ws.Cells.Item(133, 12).Value = 68550
$ws.Cells.Item(133, 13).Value = -32470
$ws.Cells.Item(133, 14).Value = -73610
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 14).ClearContents()
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()
$ws.Cells.Item(140, 8).Value = 99999.60000000001
$ws.Cells.Item(140, 10).Value = 99999.60000000001
$ws.Cells.Item(140, 12).Value = 99999.60000000001
$ws.Cells.Item(140, 14).Value = -110359.6
$ws.Cells.Item(141, 8).Value = 99693.89999999999
$ws.Cells.Item(141, 10).Value = 107548.78
$ws.Cells.Item(141, 12).Value = 107548.78
$ws.Cells.Item(141, 14).Value = -117908.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 765.25
$ws.Cells.Item(12, 9).Value = 367
$ws.Cells.Item(12, 10).Value = 898
$ws.Cells.Item(12, 11).Value = 1101
$ws.Cells.Item(12, 12).Value = 2694
$ws.Cells.Item(12, 13).Value = -928
$ws.Cells.Item(12, 14).Value = -3040
$ws.Cells.Item(46, 8).Value = 100000000
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 100000000
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 300000000
$ws.Cells.Item(46, 13).ClearContents()
$ws.Cells.Item(46, 14).Value = -300000182
$ws.Cells.Item(68, 8).Value = 2419.7097
$ws.Cells.Item(68, 10).Value = 2792.5293
$ws.Cells.Item(68, 12).Value = 8377.5879
$ws.Cells.Item(68, 14).Value = -9999.5879
$ws.Cells.Item(71, 8).Value = 2419.7097
$ws.Cells.Item(71, 10).Value = 2792.5293
$ws.Cells.Item(71, 12).Value = 25132.7637
$ws.Cells.Item(71, 14).Value = -33244.7637
$ws.Cells.Item(112, 8).Value = 11031.429
$ws.Cells.Item(112, 9).Value = 2222
$ws.Cells.Item(112, 10).Value = 12499.667
$ws.Cells.Item(112, 11).Value = 6666
$ws.Cells.Item(112, 12).Value = 37499.001
$ws.Cells.Item(112, 13).Value = -5558
$ws.Cells.Item(112, 14).Value = -39715.001
$ws.Cells.Item(113, 8).Value = 796.9167
$ws.Cells.Item(113, 9).Value = 177.77777
$ws.Cells.Item(113, 10).Value = 2654.3333
$ws.Cells.Item(113, 11).Value = 533.33331
$ws.Cells.Item(113, 12).Value = 7962.999899999999
$ws.Cells.Item(113, 13).Value = 1636.66669
$ws.Cells.Item(113, 14).Value = -12302.9999
$ws.Cells.Item(121, 8).Value = 465
$ws.Cells.Item(121, 9).Value = 465
$ws.Cells.Item(121, 11).Value = 1395
$ws.Cells.Item(121, 13).Value = -85
$ws.Cells.Item(128, 8).Value = 158015
$ws.Cells.Item(128, 9).Value = 158015
$ws.Cells.Item(128, 11).Value = 474045
$ws.Cells.Item(128, 13).Value = -469065
$ws.Cells.Item(131, 8).Value = 42679.07
$ws.Cells.Item(131, 10).Value = 10193.929
$ws.Cells.Item(131, 12).Value = 30581.787
$ws.Cells.Item(131, 14).Value = -40661.787
$ws.Cells.Item(133, 8).Value = 6679
$ws.Cells.Item(133, 9).Value = 4014.8
$ws.Cells.Item(133, 11).Value = 12044.4
$ws.Cells.Item(133, 13).Value = -6984.400000000001
$ws.Cells.Item(134, 8).Value = 2312.9583
$ws.Cells.Item(134, 9).Value = 712.5714
$ws.Cells.Item(134, 11).Value = 2137.7142
$ws.Cells.Item(134, 13).Value = 2932.2858
$ws.Cells.Item(136, 8).Value = 2492.375
$ws.Cells.Item(136, 9).Value = 2312.1667
$ws.Cells.Item(136, 11).Value = 6936.500100000001
$ws.Cells.Item(136, 13).Value = -1836.500100000001
$ws.Cells.Item(137, 8).Value = 2188.2856
$ws.Cells.Item(137, 9).Value = 1251.7778
$ws.Cells.Item(137, 11).Value = 3755.3334
$ws.Cells.Item(137, 13).Value = 1344.6666
$ws.Cells.Item(138, 8).Value = 4653
$ws.Cells.Item(138, 9).Value = 4653
$ws.Cells.Item(138, 11).Value = 13959
$ws.Cells.Item(138, 13).Value = -8819
$ws.Cells.Item(139, 8).Value = 1751.4231
$ws.Cells.Item(139, 9).Value = 1582.48
$ws.Cells.Item(139, 10).Value = 5975
$ws.Cells.Item(139, 11).Value = 4747.440000000001
$ws.Cells.Item(139, 12).Value = 17925
$ws.Cells.Item(139, 13).Value = 392.5599999999995
$ws.Cells.Item(139, 14).Value = -28205
$ws.Cells.Item(140, 8).Value = 1235.5
$ws.Cells.Item(140, 9).Value = 1235.5
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 3706.5
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 13).Value = 1473.5
$ws.Cells.Item(140, 14).ClearContents()
$ws.Cells.Item(141, 8).Value = 11321.5
$ws.Cells.Item(141, 9).Value = 2948.5
$ws.Cells.Item(141, 10).Value = 15508
$ws.Cells.Item(141, 11).Value = 8845.5
$ws.Cells.Item(141, 12).Value = 46524
$ws.Cells.Item(141, 13).Value = -3665.5
$ws.Cells.Item(141, 14).Value = -56884

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 92481.625
$ws.Cells.Item(70, 10).Value = 10066.5
$ws.Cells.Item(70, 12).Value = 10066.5
$ws.Cells.Item(70, 14).Value = -10606.5
$ws.Cells.Item(73, 8).Value = 92481.625
$ws.Cells.Item(73, 10).Value = 10066.5
$ws.Cells.Item(73, 12).Value = 10066.5
$ws.Cells.Item(73, 14).Value = -11938.5
$ws.Cells.Item(80, 8).Value = 18433.223
$ws.Cells.Item(80, 9).Value = 17985.572
$ws.Cells.Item(80, 10).Value = 20000
$ws.Cells.Item(80, 11).Value = 17985.572
$ws.Cells.Item(80, 12).Value = 20000
$ws.Cells.Item(80, 13).Value = -16987.572
$ws.Cells.Item(80, 14).Value = -21996
$ws.Cells.Item(83, 8).Value = 18433.223
$ws.Cells.Item(83, 9).Value = 17985.572
$ws.Cells.Item(83, 10).Value = 20000
$ws.Cells.Item(83, 11).Value = 89927.86
$ws.Cells.Item(83, 12).Value = 100000
$ws.Cells.Item(83, 13).Value = -84935.86
$ws.Cells.Item(83, 14).Value = -109984
$ws.Cells.Item(98, 8).Value = 37997
$ws.Cells.Item(98, 10).Value = 37997
$ws.Cells.Item(98, 12).Value = 37997
$ws.Cells.Item(98, 14).Value = -43987
$ws.Cells.Item(113, 8).Value = 3585.7856
$ws.Cells.Item(113, 9).Value = 1766.8334
$ws.Cells.Item(113, 11).Value = 1766.8334
$ws.Cells.Item(113, 13).Value = 403.1666
$ws.Cells.Item(132, 8).Value = 15037.833
$ws.Cells.Item(132, 9).Value = 12045.4
$ws.Cells.Item(132, 11).Value = 36136.2
$ws.Cells.Item(132, 13).Value = -33606.2
$ws.Cells.Item(133, 8).Value = 148999
$ws.Cells.Item(133, 10).Value = 148999
$ws.Cells.Item(133, 12).Value = 148999
$ws.Cells.Item(133, 14).Value = -159119
$ws.Cells.Item(135, 8).Value = 162222
$ws.Cells.Item(135, 9).Value = 162222
$ws.Cells.Item(135, 11).Value = 162222
$ws.Cells.Item(135, 13).Value = -157152
$ws.Cells.Item(139, 8).Value = 100000
$ws.Cells.Item(139, 10).Value = 100000
$ws.Cells.Item(139, 12).Value = 100000
$ws.Cells.Item(139, 14).Value = -110280

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7020.52
$ws.Cells.Item(7, 9).Value = 5171.3
$ws.Cells.Item(7, 11).Value = 5171.3
$ws.Cells.Item(7, 13).Value = -5059.3
$ws.Cells.Item(22, 8).Value = 3989.3333
$ws.Cells.Item(22, 9).Value = 4168.4287
$ws.Cells.Item(22, 10).Value = 3738.6
$ws.Cells.Item(22, 11).Value = 4168.4287
$ws.Cells.Item(22, 12).Value = 3738.6
$ws.Cells.Item(22, 13).Value = -3873.4287
$ws.Cells.Item(22, 14).Value = -4328.6
$ws.Cells.Item(27, 8).Value = 3989.3333
$ws.Cells.Item(27, 9).Value = 4168.4287
$ws.Cells.Item(27, 10).Value = 3738.6
$ws.Cells.Item(27, 11).Value = 4168.4287
$ws.Cells.Item(27, 12).Value = 3738.6
$ws.Cells.Item(27, 13).Value = -4061.4287
$ws.Cells.Item(27, 14).Value = -3952.6
$ws.Cells.Item(46, 8).Value = 2031.7742
$ws.Cells.Item(68, 8).Value = 3578
$ws.Cells.Item(68, 9).Value = 2158.3333
$ws.Cells.Item(68, 10).Value = 4997.6665
$ws.Cells.Item(68, 11).Value = 2158.3333
$ws.Cells.Item(68, 12).Value = 4997.6665
$ws.Cells.Item(68, 13).Value = -1409.3333
$ws.Cells.Item(68, 14).Value = -6495.6665
$ws.Cells.Item(69, 8).Value = 67500
$ws.Cells.Item(69, 10).Value = 66000
$ws.Cells.Item(69, 12).Value = 66000
$ws.Cells.Item(69, 14).Value = -67622
$ws.Cells.Item(71, 8).Value = 3578
$ws.Cells.Item(71, 9).Value = 2158.3333
$ws.Cells.Item(71, 10).Value = 4997.6665
$ws.Cells.Item(71, 11).Value = 10791.6665
$ws.Cells.Item(71, 12).Value = 24988.3325
$ws.Cells.Item(71, 13).Value = -7047.666499999999
$ws.Cells.Item(71, 14).Value = -32476.3325
$ws.Cells.Item(72, 8).Value = 67500
$ws.Cells.Item(72, 10).Value = 66000
$ws.Cells.Item(72, 12).Value = 198000
$ws.Cells.Item(72, 14).Value = -206112
$ws.Cells.Item(104, 8).Value = 13118.333
$ws.Cells.Item(104, 10).Value = 13118.333
$ws.Cells.Item(104, 12).Value = 13118.333
$ws.Cells.Item(104, 14).Value = -20106.333
$ws.Cells.Item(122, 8).Value = 4811350
$ws.Cells.Item(122, 9).Value = 3712.5264
$ws.Cells.Item(122, 10).Value = 17860652
$ws.Cells.Item(122, 11).Value = 11137.5792
$ws.Cells.Item(122, 12).Value = 53581956
$ws.Cells.Item(122, 13).Value = -8687.5792
$ws.Cells.Item(122, 14).Value = -53586856
$ws.Cells.Item(126, 8).Value = 7020.52
$ws.Cells.Item(126, 9).Value = 5171.3
$ws.Cells.Item(126, 11).Value = 15513.9
$ws.Cells.Item(126, 13).Value = -13043.9
$ws.Cells.Item(132, 8).Value = 4047.5813
$ws.Cells.Item(132, 9).Value = 3522.5312
$ws.Cells.Item(132, 11).Value = 10567.5936
$ws.Cells.Item(132, 13).Value = -8037.5936
$ws.Cells.Item(133, 8).Value = 78161.664
$ws.Cells.Item(133, 10).Value = 78161.664
$ws.Cells.Item(133, 12).Value = 78161.664
$ws.Cells.Item(133, 14).Value = -83221.664
$ws.Cells.Item(134, 8).Value = 85107.25
$ws.Cells.Item(134, 10).Value = 85107.25
$ws.Cells.Item(134, 12).Value = 85107.25
$ws.Cells.Item(134, 14).Value = -95247.25
$ws.Cells.Item(137, 8).Value = 49999.668
$ws.Cells.Item(137, 10).Value = 59999.5
$ws.Cells.Item(137, 12).Value = 59999.5
$ws.Cells.Item(137, 14).Value = -70199.5
$ws.Cells.Item(140, 8).Value = 69033.75
$ws.Cells.Item(140, 10).Value = 62692.5
$ws.Cells.Item(140, 12).Value = 62692.5
$ws.Cells.Item(140, 14).Value = -73052.5
$ws.Cells.Item(141, 8).Value = 73715
$ws.Cells.Item(141, 10).Value = 73715
$ws.Cells.Item(141, 12).Value = 73715
$ws.Cells.Item(141, 14).Value = -84075

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 14).ClearContents()
$ws.Cells.Item(88, 8).Value = 45000
$ws.Cells.Item(88, 9).Value = 45000
$ws.Cells.Item(88, 10).Value = 45000
$ws.Cells.Item(88, 11).Value = 45000
$ws.Cells.Item(88, 12).Value = 45000
$ws.Cells.Item(88, 13).Value = -44594
$ws.Cells.Item(88, 14).Value = -45812
$ws.Cells.Item(91, 8).Value = 45000
$ws.Cells.Item(91, 9).Value = 45000
$ws.Cells.Item(91, 10).Value = 45000
$ws.Cells.Item(91, 11).Value = 45000
$ws.Cells.Item(91, 12).Value = 45000
$ws.Cells.Item(91, 13).Value = -43596
$ws.Cells.Item(91, 14).Value = -47808
$ws.Cells.Item(132, 8).Value = 1793.1842
$ws.Cells.Item(132, 9).Value = 1522.8387
$ws.Cells.Item(132, 11).Value = 4568.5161
$ws.Cells.Item(132, 13).Value = -2038.5161
$ws.Cells.Item(137, 8).Value = 99178.5
$ws.Cells.Item(137, 10).Value = 99178.5
$ws.Cells.Item(137, 12).Value = 99178.5
$ws.Cells.Item(137, 14).Value = -109378.5
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 13).ClearContents()

